# Weekly update: insert two new rows of "Poroto granado" price data
# (Vega Modelo de Temuco, La Araucanía) at the top of the data block,
# pushing the previously-newest rows (which were rows 38-53) down to
# rows 40-55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38; this shifts old rows 38-53 down to 40-55
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Range("A38:R39").EntireRow.Insert()

# --- New row 38 ---
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44574
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112030
$ws.Range("G38").Value = "Poroto granado"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 210
$ws.Range("K38").Value = 28000
$ws.Range("L38").Value = 28000
$ws.Range("M38").Value = 28000
$ws.Range("N38").Value = "`$/saco 25 kilos"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 1120
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"

# --- New row 39 ---
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44574
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112030
$ws.Range("G39").Value = "Poroto granado"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Segunda"
$ws.Range("J39").Value = 85
$ws.Range("K39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = 25000
$ws.Range("N39").Value = "`$/saco 25 kilos"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1000
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
